# Auto-generated edit script applying cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.583.04"
$ws.Range("E2").Value = "  +1.89%  "
$ws.Range("D3").Value = "3.642.66"
$ws.Range("E3").Value = "  +0.25%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.43%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.62"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "195.69"
$ws.Range("E6").Value = "  +4.06%  "
$ws.Range("D7").Value = "3.636.94"
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.624"
$ws.Range("E8").Value = "  +1.68%  "
$ws.Range("E9").Value = "  -0.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.684"
$ws.Range("E10").Value = "  +0.06%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.157"
$ws.Range("E11").Value = "  +5.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.99"
$ws.Range("E12").Value = "  +1.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000298"
$ws.Range("E13").Value = "  +15.80%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.16"
$ws.Range("E14").Value = "  +1.09%  "
$ws.Range("D15").Value = "4.213.74"
$ws.Range("E15").Value = "  -0.23%  "
$ws.Range("D16").Value = "3.636.20"
$ws.Range("E16").Value = "  -0.12%  "
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.64"
$ws.Range("E18").Value = "  +2.30%  "
$ws.Range("D19").Value = "68.439.30"
$ws.Range("E19").Value = "  +1.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.67"
$ws.Range("E20").Value = "  +0.69%  "
$ws.Range("E21").Value = "  -0.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "406.39"
$ws.Range("E22").Value = "  +1.55%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.82"
$ws.Range("E23").Value = "  +20.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.28"
$ws.Range("E24").Value = "  -1.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.57"
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.97"
$ws.Range("E26").Value = "  +2.88%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.71"
$ws.Range("E27").Value = "  +3.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.92"
$ws.Range("E28").Value = "  +8.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.17"
$ws.Range("E29").Value = "  +1.87%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.14"
$ws.Range("E30").Value = "  +17.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.25"
$ws.Range("E31").Value = "  +1.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.95"
$ws.Range("E32").Value = "  +1.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "693.30"
$ws.Range("E33").Value = "  +18.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "12.31"
$ws.Range("E34").Value = "  +2.48%  "
$ws.Range("E35").Value = "  +4.71%  "
$ws.Range("E36").Value = "  -2.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "42.94"
$ws.Range("E37").Value = "  -0.52%  "
$ws.Range("E38").Value = "  +9.98%  "
$ws.Range("E39").Value = "  +10.80%  "
$ws.Range("E40").Value = "  -0.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.86"
$ws.Range("E41").Value = "  +17.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.15"
$ws.Range("E42").Value = "  +11.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.136"
$ws.Range("E43").Value = "  +1.20%  "
$ws.Range("D44").Value = "3.157.63"
$ws.Range("E44").Value = "  +16.42%  "
$ws.Range("E45").Value = "  -0.40%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0429"
$ws.Range("E46").Value = "  +3.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.78"
$ws.Range("E47").Value = "  +17.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.97"
$ws.Range("E48").Value = "  +4.72%  "
$ws.Range("E49").Value = "  -0.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "143.45"
$ws.Range("E50").Value = "  +1.46%  "
$ws.Range("B51").Value = "ApeXProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.10"
$ws.Range("E51").Value = "  -1.74%  "
